$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 203202.6
$ws.Range("I6").Value = 304375.5
$ws.Range("J6").Value = 856.8
$ws.Range("K6").Value = 913126.5
$ws.Range("L6").Value = 2570.4
$ws.Range("M6").Value = -913014.5
$ws.Range("N6").Value = -2794.4
$ws.Range("H62").Value = 1837.1428
$ws.Range("I62").Value = 1643.3334
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 1643.3334
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1019.3334
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 1837.1428
$ws.Range("I65").Value = 1643.3334
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 8216.666999999999
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -5096.666999999999
$ws.Range("N65").Value = -21240
$ws.Range("H86").Value = 6812
$ws.Range("I86").Value = 1483
$ws.Range("J86").Value = 9095.857
$ws.Range("K86").Value = 1483
$ws.Range("L86").Value = 9095.857
$ws.Range("M86").Value = -360
$ws.Range("N86").Value = -11341.857
$ws.Range("H89").Value = 6812
$ws.Range("I89").Value = 1483
$ws.Range("J89").Value = 9095.857
$ws.Range("K89").Value = 7415
$ws.Range("L89").Value = 45479.285
$ws.Range("M89").Value = -1799
$ws.Range("N89").Value = -56711.285
$ws.Range("H100").Value = 1544.2
$ws.Range("I100").Value = 1776.25
$ws.Range("J100").Value = 616
$ws.Range("K100").Value = 1776.25
$ws.Range("L100").Value = 616
$ws.Range("M100").Value = -1235.25
$ws.Range("N100").Value = -1698
$ws.Range("H106").Value = 3138.2144
$ws.Range("I106").Value = 3242.3
$ws.Range("K106").Value = 3242.3
$ws.Range("M106").Value = -2611.3
$ws.Range("H111").Value = 10800.846
$ws.Range("I111").Value = 21775.4
$ws.Range("J111").Value = 3941.75
$ws.Range("K111").Value = 65326.2
$ws.Range("L111").Value = 11825.25
$ws.Range("M111").Value = -62259.2
$ws.Range("N111").Value = -17959.25
$ws.Range("H113").Value = 101759.5
$ws.Range("I113").Value = 201619
$ws.Range("K113").Value = 201619
$ws.Range("M113").Value = -198365
$ws.Range("H132").Value = 4390511
$ws.Range("I132").Value = 5004561
$ws.Range("J132").Value = 4438.7144
$ws.Range("K132").Value = 15013683
$ws.Range("L132").Value = 13316.1432
$ws.Range("M132").Value = -15011153
$ws.Range("N132").Value = -18376.1432
$ws.Range("H134").Value = 59989.5
$ws.Range("J134").Value = 59989.5
$ws.Range("L134").Value = 59989.5
$ws.Range("N134").Value = -70129.5
$ws.Range("H141").Value = 1704.0526
$ws.Range("I141").Value = 1548.7222
$ws.Range("K141").Value = 4646.1666
$ws.Range("M141").Value = 533.8334000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 57556.555
$ws.Range("I2").Value = 2060.4167
$ws.Range("J2").Value = 168548.83
$ws.Range("K2").Value = 2060.4167
$ws.Range("L2").Value = 168548.83
$ws.Range("M2").Value = -1947.4167
$ws.Range("N2").Value = -168774.83
$ws.Range("H32").Value = 21836.576
$ws.Range("I32").Value = 3882.8289
$ws.Range("J32").Value = 173446
$ws.Range("K32").Value = 3882.8289
$ws.Range("L32").Value = 173446
$ws.Range("M32").Value = -3595.8289
$ws.Range("N32").Value = -174020
$ws.Range("H45").Value = 43927.875
$ws.Range("I45").Value = 68578.8
$ws.Range("J45").Value = 2843
$ws.Range("K45").Value = 68578.8
$ws.Range("L45").Value = 2843
$ws.Range("M45").Value = -68201.8
$ws.Range("N45").Value = -3597
$ws.Range("H88").Value = 3480.5
$ws.Range("I88").Value = 4601.25
$ws.Range("J88").Value = 2733.3333
$ws.Range("K88").Value = 4601.25
$ws.Range("L88").Value = 2733.3333
$ws.Range("M88").Value = -4195.25
$ws.Range("N88").Value = -3545.3333
$ws.Range("H91").Value = 3480.5
$ws.Range("I91").Value = 4601.25
$ws.Range("J91").Value = 2733.3333
$ws.Range("K91").Value = 4601.25
$ws.Range("L91").Value = 2733.3333
$ws.Range("M91").Value = -3197.25
$ws.Range("N91").Value = -5541.3333
$ws.Range("H110").Value = 62625692
$ws.Range("I110").Value = 66800710
$ws.Range("J110").Value = 410
$ws.Range("K110").Value = 66800710
$ws.Range("L110").Value = 410
$ws.Range("M110").Value = -66798665
$ws.Range("N110").Value = -4500
$ws.Range("H116").Value = 57556.555
$ws.Range("I116").Value = 2060.4167
$ws.Range("J116").Value = 168548.83
$ws.Range("K116").Value = 2060.4167
$ws.Range("L116").Value = 168548.83
$ws.Range("M116").Value = 233.5832999999998
$ws.Range("N116").Value = -173136.83
$ws.Range("H130").Value = 25395
$ws.Range("J130").Value = 25395
$ws.Range("L130").Value = 25395
$ws.Range("N130").Value = -35435

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 57556.555
$ws.Range("I3").Value = 2060.4167
$ws.Range("J3").Value = 168548.83
$ws.Range("K3").Value = 2060.4167
$ws.Range("L3").Value = 168548.83
$ws.Range("M3").Value = -1946.4167
$ws.Range("N3").Value = -168776.83
$ws.Range("H99").Value = 1919.9166
$ws.Range("I99").Value = 1807.3334
$ws.Range("J99").Value = 1957.4445
$ws.Range("K99").Value = 1807.3334
$ws.Range("L99").Value = 1957.4445
$ws.Range("M99").Value = -309.3334
$ws.Range("N99").Value = -4953.4445
$ws.Range("H105").Value = 155683.77
$ws.Range("I105").Value = 252145
$ws.Range("J105").Value = 112812.11
$ws.Range("K105").Value = 252145
$ws.Range("L105").Value = 112812.11
$ws.Range("M105").Value = -250398
$ws.Range("N105").Value = -116306.11

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 7010.3335
$ws.Range("I33").Value = 7010.3335
$ws.Range("K33").Value = 7010.3335
$ws.Range("M33").Value = -6631.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1287.625
$ws.Range("J17").Value = 3002
$ws.Range("L17").Value = 9006
$ws.Range("N17").Value = -9344
$ws.Range("H92").Value = 500
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H131").Value = 803.62
$ws.Range("J131").Value = 847.0345
$ws.Range("L131").Value = 2541.1035
$ws.Range("N131").Value = -12621.1035

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 1766.6666
$ws.Range("I36").Value = 1766.6666
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1766.6666
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1281.6666
$ws.Range("N36").ClearContents()
$ws.Range("H46").Value = 10757.8
$ws.Range("I46").Value = 5990
$ws.Range("J46").Value = 11949.75
$ws.Range("K46").Value = 5990
$ws.Range("L46").Value = 11949.75
$ws.Range("M46").Value = -5834
$ws.Range("N46").Value = -12261.75
$ws.Range("H113").Value = 3139.7144
$ws.Range("I113").Value = 4770.6665
$ws.Range("J113").Value = 1916.5
$ws.Range("K113").Value = 4770.6665
$ws.Range("L113").Value = 1916.5
$ws.Range("M113").Value = -2600.6665
$ws.Range("N113").Value = -6256.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 73450.64
$ws.Range("I40").Value = 334333.34
$ws.Range("J40").Value = 2300.818
$ws.Range("K40").Value = 334333.34
$ws.Range("L40").Value = 2300.818
$ws.Range("M40").Value = -334197.34
$ws.Range("N40").Value = -2572.818
$ws.Range("H46").Value = 633384.75
$ws.Range("I46").Value = 484.85715
$ws.Range("J46").Value = 1125640.2
$ws.Range("K46").Value = 484.85715
$ws.Range("L46").Value = 1125640.2
$ws.Range("M46").Value = -296.85715
$ws.Range("N46").Value = -1126016.2
$ws.Range("H122").Value = 2522.3635
$ws.Range("I122").Value = 2574.6
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7723.799999999999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -5273.799999999999
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 100001330
$ws.Range("I96").Value = 200001470
$ws.Range("J96").Value = 1180.2
$ws.Range("K96").Value = 200001470
$ws.Range("L96").Value = 1180.2
$ws.Range("M96").Value = -200000097
$ws.Range("N96").Value = -3926.2
$ws.Range("H132").Value = 1975.8679
$ws.Range("I132").Value = 2008.7858
$ws.Range("K132").Value = 6026.357400000001
$ws.Range("M132").Value = -3496.357400000001

